$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: change content from the combined "Login, Header and Footer"
#     entry to just "Login", and its status flips from "In Progress" to
#     "Pending" (copy formatting+value from B2, which is a Pending cell).
$ws.Range("A7").Value = "There should be a template for Login"
$ws.Range("B2").Copy($ws.Range("B7"))

# --- Insert five new rows (10-14) right after the existing "Navigation"
#     row (row 9), before doing that grab the two style "templates" we
#     need: a Pending row (A8:B8) and an In-Progress row (A5:B5).
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

# Row 10: Entity List - Pending
$ws.Range("A8:B8").Copy($ws.Range("A10:B10"))
$ws.Range("A10").Value = "There should be a template for Entity List"

# Row 11: Entity Form - Pending
$ws.Range("A8:B8").Copy($ws.Range("A11:B11"))
$ws.Range("A11").Value = "There should be a template for Entity Form"

# Row 12: Toolbar - In Progress
$ws.Range("A5:B5").Copy($ws.Range("A12:B12"))
$ws.Range("A12").Value = "There should be a template for Toolbar"

# Row 13: Footer - In Progress
$ws.Range("A5:B5").Copy($ws.Range("A13:B13"))
$ws.Range("A13").Value = "There should be a template for Footer"

# Row 14: Header - In Progress
$ws.Range("A5:B5").Copy($ws.Range("A14:B14"))
$ws.Range("A14").Value = "There should be a template for  Header"

# --- Row heights / borders to mirror the final layout.
$ws.Rows.Item(6).RowHeight = 15.75
$ws.Rows.Item(7).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 15.75
$ws.Rows.Item(12).RowHeight = 16.5
$ws.Rows.Item(13).RowHeight = 16.5
$ws.Rows.Item(14).RowHeight = 16.5

# --- Trailing empty spacer row.
$ws.Rows.Item(15).RowHeight = 15.75

# --- Selection / active cell.
$ws.Range("B2").Select() | Out-Null
